$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 266327.7
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 270991.3
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 812973.8999999999
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -813309.8999999999
$ws.Range("H106").Value = 1915.6666
$ws.Range("I106").Value = 1698.8
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1698.8
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1067.8
$ws.Range("N106").Value = -4262
$ws.Range("H129").Value = 935.6094000000001
$ws.Range("I129").Value = 349.9
$ws.Range("J129").Value = 1044.0741
$ws.Range("K129").Value = 1049.7
$ws.Range("L129").Value = 3132.2223
$ws.Range("M129").Value = 3950.3
$ws.Range("N129").Value = -13132.2223
$ws.Range("H132").Value = 2274547
$ws.Range("I132").Value = 1663.8918
$ws.Range("J132").Value = 14288357
$ws.Range("K132").Value = 4991.6754
$ws.Range("L132").Value = 42865071
$ws.Range("M132").Value = -2461.6754
$ws.Range("N132").Value = -42870131
$ws.Range("H138").Value = 3241.76
$ws.Range("I138").Value = 1158.8235
$ws.Range("J138").Value = 4314.788
$ws.Range("K138").Value = 3476.4705
$ws.Range("L138").Value = 12944.364
$ws.Range("M138").Value = 1663.5295
$ws.Range("N138").Value = -23224.364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7822.52
$ws.Range("I32").Value = 7188.165
$ws.Range("J32").Value = 28333.334
$ws.Range("K32").Value = 7188.165
$ws.Range("L32").Value = 28333.334
$ws.Range("M32").Value = -6901.165
$ws.Range("N32").Value = -28907.334
$ws.Range("H44").Value = 22299
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 22299
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 22299
$ws.Range("N44").Value = -23275
$ws.Range("H55").Value = 19242.4
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 19242.4
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 19242.4
$ws.Range("N55").Value = -19872.4
$ws.Range("H61").Value = 2748.8235
$ws.Range("I61").Value = 2884.182
$ws.Range("J61").Value = 2500.6667
$ws.Range("K61").Value = 2884.182
$ws.Range("L61").Value = 2500.6667
$ws.Range("M61").Value = -2672.182
$ws.Range("N61").Value = -2924.6667
$ws.Range("H63").Value = 2127.3076
$ws.Range("I63").Value = 2096.25
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2096.25
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1410.25
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2127.3076
$ws.Range("I66").Value = 2096.25
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 10481.25
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -7049.25
$ws.Range("N66").Value = -19364
$ws.Range("H80").Value = 26871.428
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 26871.428
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -28867.428
$ws.Range("H83").Value = 26871.428
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 26871.428
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -90598.284
$ws.Range("H132").Value = 1696.74
$ws.Range("I132").Value = 1424.3235
$ws.Range("J132").Value = 2275.625
$ws.Range("K132").Value = 4272.970499999999
$ws.Range("L132").Value = 6826.875
$ws.Range("M132").Value = -1742.970499999999
$ws.Range("N132").Value = -11886.875
$ws.Range("H136").Value = 2748.8235
$ws.Range("I136").Value = 2884.182
$ws.Range("J136").Value = 2500.6667
$ws.Range("K136").Value = 8652.545999999998
$ws.Range("L136").Value = 7502.000100000001
$ws.Range("M136").Value = -6102.545999999998
$ws.Range("N136").Value = -12602.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H82").Value = 10431.4
$ws.Range("I82").Value = 3539.25
$ws.Range("J82").Value = 38000
$ws.Range("K82").Value = 3539.25
$ws.Range("L82").Value = 38000
$ws.Range("M82").Value = -3156.25
$ws.Range("N82").Value = -38766
$ws.Range("H85").Value = 10431.4
$ws.Range("I85").Value = 3539.25
$ws.Range("J85").Value = 38000
$ws.Range("K85").Value = 3539.25
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = -2213.25
$ws.Range("N85").Value = -40652
$ws.Range("H94").Value = 1333.6875
$ws.Range("I94").Value = 1314.1154
$ws.Range("J94").Value = 1418.5
$ws.Range("K94").Value = 1314.1154
$ws.Range("L94").Value = 1418.5
$ws.Range("M94").Value = -863.1153999999999
$ws.Range("N94").Value = -2320.5
$ws.Range("H132").Value = 374900
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 374900
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 374900
$ws.Range("N132").Value = -385020
$ws.Range("H134").Value = 25937.666
$ws.Range("I134").Value = 2088.2222
$ws.Range("J134").Value = 68866.664
$ws.Range("K134").Value = 6264.6666
$ws.Range("L134").Value = 206599.992
$ws.Range("M134").Value = -3729.6666
$ws.Range("N134").Value = -211669.992
$ws.Range("H135").Value = 48086.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 48086.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48086.25
$ws.Range("N135").Value = -58226.25
$ws.Range("H138").Value = 97926.664
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 97926.664
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 97926.664
$ws.Range("N138").Value = -108206.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 31397.092
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 31397.092
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 31397.092
$ws.Range("N9").Value = -31733.092
$ws.Range("H87").Value = 49800
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 49800
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 49800
$ws.Range("N87").Value = -52172
$ws.Range("H90").Value = 49800
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 49800
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 149400
$ws.Range("N90").Value = -161256
$ws.Range("H99").Value = 1959.0769
$ws.Range("I99").Value = 1953.7142
$ws.Range("J99").Value = 1972.7273
$ws.Range("K99").Value = 1953.7142
$ws.Range("L99").Value = 1972.7273
$ws.Range("M99").Value = -455.7141999999999
$ws.Range("N99").Value = -4968.7273
$ws.Range("H126").Value = 1959.0769
$ws.Range("I126").Value = 1953.7142
$ws.Range("J126").Value = 1972.7273
$ws.Range("K126").Value = 5861.142599999999
$ws.Range("L126").Value = 5918.1819
$ws.Range("M126").Value = -3391.142599999999
$ws.Range("N126").Value = -10858.1819
$ws.Range("H132").Value = 1767.7
$ws.Range("I132").Value = 1319.4546
$ws.Range("J132").Value = 2315.5557
$ws.Range("K132").Value = 3958.3638
$ws.Range("L132").Value = 6946.6671
$ws.Range("M132").Value = -1428.3638
$ws.Range("N132").Value = -12006.6671
$ws.Range("H134").Value = 41667870
$ws.Range("I134").Value = 1427.875
$ws.Range("J134").Value = 125000760
$ws.Range("K134").Value = 4283.625
$ws.Range("L134").Value = 375002280
$ws.Range("M134").Value = -1748.625
$ws.Range("N134").Value = -375007350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 70749.836
$ws.Range("I70").Value = 202245
$ws.Range("J70").Value = 5002.25
$ws.Range("K70").Value = 202245
$ws.Range("L70").Value = 5002.25
$ws.Range("M70").Value = -201975
$ws.Range("N70").Value = -5542.25
$ws.Range("H73").Value = 70749.836
$ws.Range("I73").Value = 202245
$ws.Range("J73").Value = 5002.25
$ws.Range("K73").Value = 202245
$ws.Range("L73").Value = 5002.25
$ws.Range("M73").Value = -201309
$ws.Range("N73").Value = -6874.25
$ws.Range("H80").Value = 3732.75
$ws.Range("I80").Value = 3890.3333
$ws.Range("J80").Value = 3260
$ws.Range("K80").Value = 3890.3333
$ws.Range("L80").Value = 3260
$ws.Range("M80").Value = -2892.3333
$ws.Range("N80").Value = -5256
$ws.Range("H83").Value = 3732.75
$ws.Range("I83").Value = 3890.3333
$ws.Range("J83").Value = 3260
$ws.Range("K83").Value = 19451.6665
$ws.Range("L83").Value = 16300
$ws.Range("M83").Value = -14459.6665
$ws.Range("N83").Value = -26284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 15185
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 15185
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 15185
$ws.Range("N104").Value = -22173
$ws.Range("H136").Value = 3360.0166
$ws.Range("I136").Value = 1777.1163
$ws.Range("J136").Value = 7363.8237
$ws.Range("K136").Value = 5331.3489
$ws.Range("L136").Value = 22091.4711
$ws.Range("M136").Value = -2781.3489
$ws.Range("N136").Value = -27191.4711
